# Auto-generated Excel COM-interop script to apply cryptos.xlsx diff
# Commit: Updated cryptos list on Wed Oct 11 19:39:10 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'26.743.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.47%  "
$ws.Range("D3").Value = "'1.561.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'206.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.488"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'21.90"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").Value = "'0.0583"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("D11").Value = "'0.0860"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "'1.783.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "'1.564.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").Value = "'3.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "'61.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.21%  "
$ws.Range("D17").Value = "'26.768.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'213.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'7.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("D20").Value = "'0.0₃0675"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("E23").Value = "  -2.08%  "
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").Value = "'152.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("D26").Value = "'6.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").Value = "'14.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("D32").Value = "'3.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("D33").Value = "'1.384.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("D34").Value = "'2.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("D37").Value = "'0.930"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.50%  "
$ws.Range("E38").Value = "  -2.58%  "
$ws.Range("D39").Value = "'0.520"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.92%  "
$ws.Range("D40").Value = "'0.812"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("D43").Value = "'5.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("E45").Value = "  -2.17%  "
$ws.Range("D46").Value = "'63.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").Value = "'1.696.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").Value = "'85.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").Value = "'0.0₇0984"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0947"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0492"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.51%  "
